$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text cell (A1) with the new exchange rate figures ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 13.78 = 55720.74 pesos`n✅ 55720.74 pesos = 13.71 = 977.49 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet numeric cells ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 72.55
$tasas.Range("O10").Value = 4042.54
$tasas.Range("N12").Value = 4065
$tasas.Range("O12").Value = 71.311
